$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "RP Jasa" header in F1, sharing the bold style used by the
# other header cells (A1:E1) so it reuses the same cell style index.
$ws.Range("F1").Value = "RP Jasa"
$ws.Range("F1").Font.Bold = $true

# Column F gets its own (slightly different) custom width.
$ws.Columns("F").ColumnWidth = 17.25

# The saved view now has F4 as the active/selected cell.
[void]$ws.Range("F4").Select()
